$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new "through" date
$ws.Name = "Through 2022-02-22"

# Update the header label cell (I1) which holds the "2022 (through 02-21)" text
$ws.Range("I1").Value = "2022 (through 02-22)"

# Update February 2022 total (I3) and the yearly Total row (I14)
$ws.Range("I3").Value = 111
$ws.Range("I14").Value = 270
